$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.435.52'
$ws.Range('E2').Value = '  +1.69%  '
$ws.Range('D3').Value = '2.281.31'
$ws.Range('E3').Value = '  +0.77%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.27'
$ws.Range('E5').Value = '  +0.95%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.99'
$ws.Range('E6').Value = '  +6.53%  '
$ws.Range('E7').Value = '  +0.26%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.496'
$ws.Range('E9').Value = '  +3.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.08'
$ws.Range('E10').Value = '  +11.56%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0803'
$ws.Range('E11').Value = '  +0.79%  '
$ws.Range('E12').Value = '  -2.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.72'
$ws.Range('E13').Value = '  +1.10%  '
$ws.Range('D14').Value = '2.637.98'
$ws.Range('E14').Value = '  +0.91%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.50'
$ws.Range('E15').Value = '  +2.07%  '
$ws.Range('D16').Value = '2.295.03'
$ws.Range('E16').Value = '  +1.96%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.801'
$ws.Range('E17').Value = '  +4.72%  '
$ws.Range('D18').Value = '42.333.74'
$ws.Range('E18').Value = '  +1.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.62'
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('D20').Value = '0.0₃0914'
$ws.Range('E20').Value = '  +1.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.03'
$ws.Range('E21').Value = '  +1.64%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.82'
$ws.Range('E22').Value = '  +1.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '242.06'
$ws.Range('E23').Value = '  +0.93%  '
$ws.Range('E24').Value = '  +0.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.95'
$ws.Range('E25').Value = '  +1.51%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.91'
$ws.Range('E27').Value = '  -0.34%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.03'
$ws.Range('E28').Value = '  +10.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.57'
$ws.Range('E29').Value = '  +0.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.12'
$ws.Range('E30').Value = '  +2.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '160.96'
$ws.Range('E31').Value = '  +0.60%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.26'
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.13'
$ws.Range('E34').Value = '  +4.28%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0747'
$ws.Range('E35').Value = '  +0.63%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.22'
$ws.Range('E36').Value = '  +2.17%  '
$ws.Range('E37').Value = '  +2.42%  '
$ws.Range('E38').Value = '  +0.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.86'
$ws.Range('E39').Value = '  +3.38%  '
$ws.Range('E40').Value = '  -0.80%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.17'
$ws.Range('E41').Value = '  +5.75%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.40'
$ws.Range('E42').Value = '  +14.41%  '
$ws.Range('D43').Value = '2.003.34'
$ws.Range('E43').Value = '  -1.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.22'
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('E45').Value = '  +3.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.01'
$ws.Range('E46').Value = '  +4.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.08'
$ws.Range('E47').Value = '  -3.31%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '53.51'
$ws.Range('E48').Value = '  +3.08%  '
$ws.Range('E49').Value = '  +1.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '93.29'
$ws.Range('E50').Value = '  +2.33%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.42'
$ws.Range('E51').Value = '  -0.28%  '
